$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-497) holds the "Förändrad" (Changed) date.
# It moves forward by one day: 2023-09-19 (serial 45188) -> 2023-09-20 (serial 45189).
$ws.Range("C2:C497").Value = 45189
